# A new weekly price-report row (Poroto verde, Región de Arica y Parinacota,
# 2022-08-16) is inserted at the top of the data set (row 84, right after the
# frozen/unchanged rows above it). Excel's usual "insert row, shift cells
# down" behaviour pushes every existing record from row 84 downward by one
# row, so the former last row (169) becomes row 170 and the sheet's used
# range grows from A1:R169 to A1:R170.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 84; everything that was at 84..169 moves to 85..170.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 44789
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112031
$ws.Range("G84").Value = "Poroto verde"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 150
$ws.Range("K84").Value = 35000
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = 35000
$ws.Range("N84").Value = "$/malla 25 kilos"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 1400
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
